# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled run).
# Source data only touches columns D (Price) and E (Volume(1h)) for nearly every
# row; rows 50/51 (EnergySwap / Decentraland) also swapped rank order, so all
# four columns (B/C/D/E) are rewritten there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" strings (e.g. "1.0000", "0.9999") are digits-and-dots that
# Excel's Range.Value setter would otherwise auto-coerce into a number,
# dropping the trailing/leading zeros that make the original text formatting.
# Force those particular assignments to be stored as text, then drop the
# style back to Normal so we don't leave a stray "@" (Text) format applied
# to the cell (matches the source, which carries no cell style at all here).
function Set-TextValue($range, $value) {
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.Style = 'Normal'
}

# Bitcoin
$ws.Range('D2').Value = '25.976.80'
$ws.Range('E2').Value = '  -0.30%  '

# Ethereum
$ws.Range('D3').Value = '1.744.12'

# TetherUSD
$ws.Range('E4').Value = '  -0.11%  '

# BNB
Set-TextValue $ws.Range('D5') '250.98'
$ws.Range('E5').Value = '  +7.77%  '

# USDC
Set-TextValue $ws.Range('D6') '1.0000'
$ws.Range('E6').Value = '  -0.07%  '

# XRP
Set-TextValue $ws.Range('D7') '0.5144'
$ws.Range('E7').Value = '  -2.25%  '

# Cardano
Set-TextValue $ws.Range('D8') '0.2764'
$ws.Range('E8').Value = '  -0.14%  '

# Dogecoin
$ws.Range('E9').Value = '  +0.14%  '

# WrappedEther
$ws.Range('D10').Value = '1.744.03'
$ws.Range('E10').Value = '  -0.25%  '

# TRON
Set-TextValue $ws.Range('D11') '0.07236'
$ws.Range('E11').Value = '  +0.65%  '

# Solana
Set-TextValue $ws.Range('D12') '15.22'
$ws.Range('E12').Value = '  -0.48%  '

# Polygon
Set-TextValue $ws.Range('D13') '0.6506'
$ws.Range('E13').Value = '  +1.55%  '

# Polkadot
Set-TextValue $ws.Range('D14') '4.637'
$ws.Range('E14').Value = '  +1.14%  '

# Litecoin
Set-TextValue $ws.Range('D15') '77.82'
$ws.Range('E15').Value = '  -0.46%  '

# Dai
Set-TextValue $ws.Range('D16') '0.9999'
$ws.Range('E16').Value = '  -0.04%  '

# BinanceUSD
Set-TextValue $ws.Range('D17') '0.9998'
$ws.Range('E17').Value = '  -0.17%  '

# WrappedBTC
$ws.Range('D18').Value = '26.004.91'
$ws.Range('E18').Value = '  +0.06%  '

# Avalanche
$ws.Range('E19').Value = '  +2.39%  '

# ShibaInu
Set-TextValue $ws.Range('D20') '0.000006808'
$ws.Range('E20').Value = '  +1.47%  '

# WrappedliquidstakedEther2.0
$ws.Range('D21').Value = '1.966.82'
$ws.Range('E21').Value = '  -0.54%  '

# Uniswap
Set-TextValue $ws.Range('D22') '4.311'
$ws.Range('E22').Value = '  -0.14%  '

# Cosmos
Set-TextValue $ws.Range('D23') '8.692'
$ws.Range('E23').Value = '  -1.05%  '

# Chainlink
Set-TextValue $ws.Range('D24') '5.390'
$ws.Range('E24').Value = '  +3.72%  '

# Monero
Set-TextValue $ws.Range('D25') '135.82'
$ws.Range('E25').Value = '  -2.39%  '

# Toncoin
Set-TextValue $ws.Range('D26') '1.509'
$ws.Range('E26').Value = '  -0.69%  '

# EthereumClassic
$ws.Range('E27').Value = '  +0.16%  '

# LidoDAOToken
Set-TextValue $ws.Range('D28') '1.786'
$ws.Range('E28').Value = '  -1.49%  '

# BitcoinCash
Set-TextValue $ws.Range('D29') '106.14'
$ws.Range('E29').Value = '  +1.85%  '

# InternetComputer(DFINITY)
Set-TextValue $ws.Range('D30') '3.951'
$ws.Range('E30').Value = '  +4.87%  '

# Stellar
$ws.Range('E31').Value = '  -0.85%  '

# Filecoin
Set-TextValue $ws.Range('D32') '3.677'
$ws.Range('E32').Value = '  +0.81%  '

# Hedera
Set-TextValue $ws.Range('D33') '0.04671'
$ws.Range('E33').Value = '  +2.87%  '

# HuobiToken
$ws.Range('E34').Value = '  +0.47%  '

# ARBITRUM
Set-TextValue $ws.Range('D35') '1.001'
$ws.Range('E35').Value = '  +0.21%  '

# ImmutableX
Set-TextValue $ws.Range('D36') '0.6258'
$ws.Range('E36').Value = '  -0.97%  '

# MXToken
Set-TextValue $ws.Range('D37') '2.736'
$ws.Range('E37').Value = '  +1.09%  '

# VeChain
Set-TextValue $ws.Range('D38') '0.01612'
$ws.Range('E38').Value = '  +1.10%  '

# RenderToken
Set-TextValue $ws.Range('D39') '1.932'
$ws.Range('E39').Value = '  +0.05%  '

# PaxDollar
Set-TextValue $ws.Range('D40') '0.9996'
$ws.Range('E40').Value = '  -0.06%  '

# Quant
Set-TextValue $ws.Range('D41') '100.64'
$ws.Range('E41').Value = '  +2.44%  '

# TheSandbox
Set-TextValue $ws.Range('D42') '0.3888'
$ws.Range('E42').Value = '  -0.31%  '

# TrustWalletToken
Set-TextValue $ws.Range('D43') '0.7590'
$ws.Range('E43').Value = '  +3.64%  '

# FraxShare
Set-TextValue $ws.Range('D44') '5.028'
$ws.Range('E44').Value = '  -0.22%  '

# Aptos
Set-TextValue $ws.Range('D45') '6.354'
$ws.Range('E45').Value = '  +0.75%  '

# Algorand
$ws.Range('E46').Value = '  -0.58%  '

# Aave
Set-TextValue $ws.Range('D47') '55.52'
$ws.Range('E47').Value = '  +3.09%  '

# Cronos
Set-TextValue $ws.Range('D48') '0.05234'
$ws.Range('E48').Value = '  -2.08%  '

# Elrond
Set-TextValue $ws.Range('D49') '30.78'
$ws.Range('E49').Value = '  +0.82%  '

# Rows 50/51 swap: Decentraland now ranks above EnergySwap.
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D50') '0.3454'
$ws.Range('E50').Value = '  -0.14%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '7.574'
$ws.Range('E51').Value = '  -1.25%  '
